$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 48, pushing the existing row 48
# (and everything below it) down by one.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with a new data record. All
# non-date/non-volume fields mirror the template used by the rest of
# the "Femacal de La Calera - Ciboulette" rows; only the date (D) and
# volume (J) values are new for this record.
$ws.Range("A48").Value = 3
$ws.Range("B48").Value = "Femacal de La Calera"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44749
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = 100112039
$ws.Range("G48").Value = "Ciboulette"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 120
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 1500
$ws.Range("M48").Value = 1500
$ws.Range("N48").Value = "`$/docena de atados"
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 500
$ws.Range("Q48").Value = 3
$ws.Range("R48").Value = "Hortaliza"
